# fix: Remover TRATAMENTO EMAGRECIMENTO e ajustar budget saude
# Updates the "Compras" budget on the Dashboard/Categorias sheets from
# R$ 2.000 to R$ 4.000 and propagates the resulting changes to "Real"
# gasto values (4669.66 -> 2169.66) and the derived totals/percentages
# across the Dashboard, Mensal, Categorias and Dados sheets. Also
# refreshes the "last updated" timestamp stamps.

$wb = $excel.ActiveWorkbook

# ---- Dashboard sheet ----
$dash = $wb.Worksheets.Item("Dashboard")

# Header timestamp
$dash.Range("A2").Value = "Atualizado: 31/12/2025 11:23"

# RESUMO DO MES - "Gastos Variaveis" row (row 7)
$dash.Range("B7").Value = 19800
$dash.Range("C7").Value = 6765.58
# D7 stores a literal percentage string (not a numeric percent), so use
# a leading quote to force text entry and avoid Excel auto-converting
# "-65%" into the number -0.65.
$dash.Range("D7").Value = "'-65%"

# RESUMO DO MES - "Obra" row (row 9), Real column recalculated
$dash.Range("C9").Value = 38734.42

# GASTOS POR CATEGORIA - "Compras" row (row 18)
$dash.Range("B18").Value = 4000
$dash.Range("C18").Value = 2169.66
$dash.Range("D18").Value = 54

# ---- Mensal sheet ----
$mensal = $wb.Worksheets.Item("Mensal")
$mensal.Range("B8").Value = 2169.66

# ---- Categorias sheet ----
$categorias = $wb.Worksheets.Item("Categorias")
$categorias.Range("C8").Value = 2169.66
$categorias.Range("D8").Value = 1830.34
$categorias.Range("E8").Value = 0.542415

# ---- Dados sheet ----
$dados = $wb.Worksheets.Item("Dados")
$dados.Range("B3").Value = "2025-12-31T11:23:58.847300"
$dados.Range("D12").Value = 2169.66
